# feat: add 2022-Q1 data
#
# 1) Insert a new worksheet "2022-Q1" right before the "总计" (totals) sheet,
#    populated with the same column layout as the other quarterly sheets.
# 2) Prepend a "2022-Q1" row to the "总计" sheet's summary table and
#    renumber the existing index column sequentially.

$wb = $excel.ActiveWorkbook

# Remember the totals sheet by name - its numeric Index will shift once we
# insert a new sheet in front of it, so re-resolve it by name afterwards.
$totalsName = $wb.Worksheets.Item(4).Name

# ---------------------------------------------------------------------
# 1. Create the "2022-Q1" worksheet immediately before "总计"
# ---------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item($totalsName)
$ns = $wb.Worksheets.Add($beforeSheet)
$ns.Name = "2022-Q1"

# Match the page margins used throughout the rest of the workbook
# (0.75in/0.75in/1in/1in/0.5in/0.5in == 54/54/72/72/36/36 points).
$ns.PageSetup.LeftMargin = 54
$ns.PageSetup.RightMargin = 54
$ns.PageSetup.TopMargin = 72
$ns.PageSetup.BottomMargin = 72
$ns.PageSetup.HeaderMargin = 36
$ns.PageSetup.FooterMargin = 36

# Borrow the bold/centered/bordered header style already used on the other
# quarterly sheets (e.g. "2021-Q1") instead of re-building it by hand.
$styleSrc = $wb.Worksheets.Item("2021-Q1")
$styleSrc.Range("B1:H1").Copy()
$ns.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats

$ns.Cells.Item(1,2).Value = "基金代码"
$ns.Cells.Item(1,3).Value = "基金名称"
$ns.Cells.Item(1,4).Value = "基金规模"
$ns.Cells.Item(1,5).Value = "股票总仓位"
$ns.Cells.Item(1,6).Value = "仓位占比"
$ns.Cells.Item(1,7).Value = "持有市值(亿元)"
$ns.Cells.Item(1,8).Value = "仓位排名"

# Helper: write $text into (row,col) of $sheet as a genuine text cell, even
# when $text looks like a number (e.g. "003397", "0.0170") so leading /
# trailing zeros survive. Cross-sheet copy/paste is unreliable here, so the
# scratch cell lives on the same sheet, far below the real data, and is
# deleted again once used.
function Set-TextCell($sheet, $row, $col, $text) {
    $scratchRow = 100
    $sheet.Cells.Item($scratchRow, 1).Formula = '="' + $text + '"'
    $sheet.Cells.Item($scratchRow, 1).Copy()
    $sheet.Cells.Item($row, $col).PasteSpecial(-4163)   # xlPasteValues
    $sheet.Cells.Item($scratchRow, 1).ClearContents()
}

# Row 2
$ns.Cells.Item(2,1).Value = 0
Set-TextCell $ns 2 2 "003397"
Set-TextCell $ns 2 3 "银华体育文化灵活配置混合"
Set-TextCell $ns 2 4 "0.39"
Set-TextCell $ns 2 5 "83.61"
Set-TextCell $ns 2 6 "4.36"
Set-TextCell $ns 2 7 "0.0170"
$ns.Cells.Item(2,8).Value = 6

# Row 3
$ns.Cells.Item(3,1).Value = 1
Set-TextCell $ns 3 2 "003659"
Set-TextCell $ns 3 3 "山西证券策略精选灵活配置混合"
Set-TextCell $ns 3 4 "0.31"
Set-TextCell $ns 3 5 "84.52"
Set-TextCell $ns 3 6 "2.94"
Set-TextCell $ns 3 7 "0.0091"
$ns.Cells.Item(3,8).Value = 9

# Match the index-column style (centered/bold/bordered, s="2") used by the
# sibling quarterly sheets.
$styleSrc.Cells.Item(2,1).Copy()
$ns.Cells.Item(2,1).PasteSpecial(-4122)
$ns.Cells.Item(3,1).PasteSpecial(-4122)

# Drop the scratch row used by Set-TextCell.
$ns.Rows.Item(100).Delete()

# ---------------------------------------------------------------------
# 2. Add a "2022-Q1" row to the "总计" summary sheet
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item($totalsName)

$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()
$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q1"
$total.Cells.Item(2,3).Value = 2
$total.Cells.Item(2,4).Value = 0.03

# Give the new index cell (A2) the same style as the other index cells.
$total.Cells.Item(3,1).Copy()
$total.Cells.Item(2,1).PasteSpecial(-4122)

# Renumber the (pre-existing) index column sequentially: 1, 2, 3.
$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(5,1).Value = 3

# Restore "2020-Q4" as the active/selected tab (it was active before this
# edit; creating the new sheet would otherwise steal the selection).
$wb.Worksheets.Item("2020-Q4").Activate()
